$d = $word.ActiveDocument

$replacements = @(
    @(" And the earth was without form, and void; and darkness was upon the face of the deep. And the Spirit of God moved upon the face of the waters.", " And the earth was without form and void; and darkness was upon the face of the deep. And the Spirit of God moved upon the face of the waters."),
    @(" 起初，　神創造天地。", " 起初，　神 創造 天 地。"),
    @(" 地是空虛混沌，淵面黑暗；　神的靈運行在水面上。", " 地 是 空虛 混沌，淵面 黑暗；　神 的靈 運行 在 水 面上。"),
    @(" 神說：「要有光」，就有了光。", " 神 說：「要有 光」，就有了 光。"),
    @(" 神看光是好的，就把光暗分開了。", " 神 看 光 是好的， 就把 光 暗 分開了。"),
    @(" 神就照著自己的形像造人，乃是照著他的形像造男造女。", " 神 就照著自己的形像 造 人，乃是照著他 的形像 造 男 造 女。"),
    @(" 神就賜福給他們，又對他們說：「要生養眾多，遍滿地面，治理這地，也要管理海裡的魚、空中的鳥，和地上各樣行動的活物。」", " 神 就賜福 給他們， 又對他們 說：「要生養 眾多，遍滿 地面，治理這地，也要管理 海 裡的魚、空 中的鳥，和地上 各樣 行動的 活物。」"),
    @(" 神說：「看哪，我將遍地上一切結種子的菜蔬和一切樹上所結有核的果子全賜給你們作食物。", " 神 說：「看哪，我將 遍 地 上 一切 結 種子 的菜蔬 和 一切 樹 上 所 結 有核 的果子 全賜 給你們 作 食物。"),
    @(" 至於地上的走獸和空中的飛鳥，並各樣爬在地上有生命的物，我將青草賜給牠們作食物。」事就這樣成了。", " 至於 地上 的走獸和 空中 的飛鳥，並各樣 爬 在 地上有 生命 的物，我將 青 草 賜給牠們作食物。」事就這樣 成了。"),
    @(" 神看著一切所造的都甚好。有晚上，有早晨，是第六日。", " 神 看著 一切 所 造的 都甚 好。有 晚上，有 早晨，是第六 日。"),
    @(" 神的兒子，耶穌基督福音的起頭。", " 神的 兒子，耶穌 基督 福音的 起頭。"),
    @(" 正如先知以賽亞（有古卷沒有以賽亞三個字）書上記著說：看哪，我要差遣我的使者在你前面，預備道路。", " 正如 先知 以賽亞（有古卷沒有以賽亞三個字）書上記著說：看哪，我要差遣 我的 使者 在 你 前面， 預備 道路。"),
    @(" 在曠野有人聲喊著說：預備主的道，修直他的路。", " 在 曠野 有人聲 喊著說：預備 主的 道，修 直 他的 路。"),
    @(" 約翰下監以後，耶穌來到加利利，宣傳　神的福音，", " 約翰 下監 以後，耶穌 來 到 加利利，宣傳　神的 福音，"),
    @(" 說：「日期滿了，　神的國近了。你們當悔改，信福音！」", " 說：「 日期 滿了，　神的 國 近了。你們當悔改， 信 福音！」"),
    @(" 耶穌順著加利利的海邊走，看見西門和西門的兄弟安得烈在海裡撒網；他們本是打魚的。", " 耶穌順著 加利利的 海邊 走，看見 西門 和 西門的 兄弟 安得烈 在 海 裡撒網；他們本 是 打魚的。"),
    @(" 耶穌對他們說：「來跟從我，我要叫你們得人如得魚一樣。」", " 耶穌 對他們 說：「來 跟從 我， 我要叫 你們 得人 如得魚一樣。」"),
    @(" 他們就立刻捨了網，跟從了他。", " 他們就 立刻 捨了 網，跟從了 他。"),
    @(" 耶穌稍往前走，又見西庇太的兒子雅各和雅各的兄弟約翰在船上補網。", " 耶穌稍 往前走，又見 西庇太的 兒子雅各 和 雅各的 兄弟 約翰 在 船上 補 網。"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND:" $old
    }
}
